$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @{
    "J2" = 1.08
    "L2" = 1.4
    "N3" = 1.85
    "O3" = 2.05
    "G5" = 1.27
    "H5" = 5.5
    "I5" = 11
    "J5" = 1.03
    "L5" = 1.27
    "M5" = 3.4
    "N5" = 1.98
    "O5" = 1.83
    "R5" = 2.75
    "S5" = 1.4
    "T5" = 5.5
    "V5" = 10
    "W5" = 7
    "X5" = 15
    "AA5" = 11
    "AC5" = 151
    "J6" = 1.08
    "L6" = 1.58
    "J7" = 1.13
    "L7" = 1.69
    "G21" = 1.91
    "R21" = 1.83
    "S21" = 1.83
    "G22" = 2.35
    "J22" = 1.06
    "K22" = 10
    "N22" = 2.05
    "O22" = 1.75
    "R22" = 1.8
    "S22" = 1.95
    "G23" = 2.4
    "R23" = 1.73
    "G24" = 1.85
    "R24" = 1.91
    "S24" = 1.8
    "R25" = 1.53
    "S25" = 2.38
    "G27" = 4.9
    "H27" = 3.55
    "I27" = 1.6
    "N27" = 1.88
    "O27" = 1.72
    "P27" = 1.39
    "Q27" = 2.42
    "T27" = 10.25
    "U27" = 22
    "V27" = 13.5
    "W27" = 65
    "X27" = 40
    "Y27" = 45
    "Z27" = 9
    "AA27" = 6.1
    "AB27" = 14
    "AC27" = 65
    "AE27" = 5.5
    "AF27" = 6.2
    "AG27" = 6.9
    "AH27" = 9.75
    "AI27" = 11
    "AJ27" = 22
    "G28" = 3.65
    "I28" = 1.98
    "N28" = 2.12
    "T28" = 7.9
    "V28" = 10.5
    "W28" = 40
    "X28" = 28
    "Y28" = 35
    "AA28" = 5.2
    "AC28" = 60
    "AE28" = 5.2
    "AF28" = 7.3
    "AG28" = 7.3
    "AH28" = 14
    "AI28" = 14.5
    "AJ28" = 26
    "J29" = 1.02
    "K29" = 19
    "P29" = 1.25
    "P30" = 1.22
    "G31" = 1.95
    "I31" = 3.8
    "R31" = 1.75
    "S31" = 2
    "U31" = 9.5
    "W31" = 17
    "AA31" = 6.5
    "AD31" = 201
    "AE31" = 11
    "AF31" = 19
    "AI31" = 29
    "AJ31" = 34
    "H34" = 4.1
    "K34" = 17
    "L34" = 1.18
    "M34" = 4.5
    "N34" = 1.6
    "O34" = 2.3
    "P34" = 1.29
    "Q34" = 3.5
    "R34" = 1.62
    "S34" = 2.2
    "T34" = 9
    "AE34" = 15
    "AJ34" = 29
    "P39" = 1.62
    "N40" = 1.9
    "O40" = 1.9
    "P40" = 1.36
    "P42" = 1.14
    "J43" = 1.03
    "K43" = 10.5
    "N43" = 1.53
    "O43" = 2.38
    "P43" = 1.29
    "N58" = 1.9
    "O58" = 1.9
}

foreach ($addr in $edits.Keys) {
    $ws.Range($addr).Value = $edits[$addr]
}
